$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

$ws.Range("A6").Value = "Data"
$ws.Range("B6").Value = 2915004
